$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Table" component (row 12) used to carry static "columns"/"rows" JSON
# props. It now points at a remote API instead, keeping "stickyHeader" and
# dropping the two JSON blobs.
#
# Before: E12=columns F12=[json columns] G12=rows H12=[json rows] I12=stickyHeader J12=TRUE
# After:  E12=api     F12=<api url>               G12=stickyHeader         H12=TRUE

$ws.Range("E12").Value2 = "api"

$ws.Range("F12").Clear()
$ws.Range("F12").Value2 = "http://localhost:9001/tableData"

$ws.Range("G12").Clear()
$ws.Range("G12").Value2 = "stickyHeader"

$ws.Range("H12").Clear()
$ws.Range("I12").Clear()

$ws.Range("H12").Value2 = $true
$ws.Range("J12").Clear()
